$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 41667090
$ws.Range("I33").Value = 52632084
$ws.Range("K33").Value = 52632084
$ws.Range("M33").Value = -52631855
$ws.Range("H101").Value = 929.7059
$ws.Range("I101").Value = 692.0833
$ws.Range("J101").Value = 1500
$ws.Range("K101").Value = 2076.2499
$ws.Range("L101").Value = 4500
$ws.Range("M101").Value = -454.2498999999998
$ws.Range("N101").Value = -7744
$ws.Range("H115").Value = 1742.5
$ws.Range("H116").Value = 2540
$ws.Range("I116").Value = 2425
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2425
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1017
$ws.Range("N116").Value = -9884
$ws.Range("H137").Value = 10001422
$ws.Range("I137").Value = 1024.6
$ws.Range("K137").Value = 3073.8
$ws.Range("M137").Value = -523.7999999999997
$ws.Range("H138").Value = 1974.14
$ws.Range("I138").Value = 1317.2084
$ws.Range("J138").Value = 2580.5386
$ws.Range("K138").Value = 3951.6252
$ws.Range("L138").Value = 7741.6158
$ws.Range("M138").Value = 1188.3748
$ws.Range("N138").Value = -18021.6158
$ws.Range("H139").Value = 39933.332
$ws.Range("J139").Value = 39933.332
$ws.Range("L139").Value = 39933.332
$ws.Range("N139").Value = -50213.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 49016.25
$ws.Range("J24").Value = 49016.25
$ws.Range("L24").Value = 49016.25
$ws.Range("N24").Value = -49764.25
$ws.Range("H32").Value = 10063.96
$ws.Range("I32").Value = 11518
$ws.Range("J32").Value = 4247.8
$ws.Range("K32").Value = 11518
$ws.Range("L32").Value = 4247.8
$ws.Range("M32").Value = -11231
$ws.Range("N32").Value = -4821.8
$ws.Range("H100").Value = 49016.25
$ws.Range("J100").Value = 49016.25
$ws.Range("L100").Value = 49016.25
$ws.Range("N100").Value = -51180.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2389.6843
$ws.Range("I107").Value = 2471.3333
$ws.Range("J107").Value = 2083.5
$ws.Range("K107").Value = 2471.3333
$ws.Range("L107").Value = 2083.5
$ws.Range("M107").Value = -551.3332999999998
$ws.Range("N107").Value = -5923.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1096.5264
$ws.Range("I16").Value = 896.9167
$ws.Range("J16").Value = 1438.7142
$ws.Range("K16").Value = 896.9167
$ws.Range("L16").Value = 1438.7142
$ws.Range("M16").Value = -609.9167
$ws.Range("N16").Value = -2012.7142
$ws.Range("H107").Value = 764.8889
$ws.Range("I107").Value = 754.2
$ws.Range("J107").Value = 778.25
$ws.Range("K107").Value = 754.2
$ws.Range("L107").Value = 778.25
$ws.Range("M107").Value = 1165.8
$ws.Range("N107").Value = -4618.25
$ws.Range("H113").Value = 1096.5264
$ws.Range("I113").Value = 896.9167
$ws.Range("J113").Value = 1438.7142
$ws.Range("K113").Value = 896.9167
$ws.Range("L113").Value = 1438.7142
$ws.Range("M113").Value = 1273.0833
$ws.Range("N113").Value = -5778.7142
$ws.Range("H118").Value = 38742
$ws.Range("J118").Value = 38742
$ws.Range("L118").Value = 38742
$ws.Range("N118").Value = -42056
$ws.Range("H122").Value = 1258.1389
$ws.Range("I122").Value = 1276.24
$ws.Range("J122").Value = 1217
$ws.Range("K122").Value = 3828.72
$ws.Range("L122").Value = 3651
$ws.Range("M122").Value = -1378.72
$ws.Range("N122").Value = -8551
$ws.Range("H140").Value = 27650
$ws.Range("J140").Value = 27650
$ws.Range("L140").Value = 27650
$ws.Range("N140").Value = -38010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3109.875
$ws.Range("I5").Value = 434.8
$ws.Range("K5").Value = 1304.4
$ws.Range("M5").Value = -1192.4
$ws.Range("H6").Value = 109.181816
$ws.Range("I6").Value = 60.142857
$ws.Range("J6").Value = 195
$ws.Range("K6").Value = 180.428571
$ws.Range("L6").Value = 585
$ws.Range("M6").Value = -67.42857100000001
$ws.Range("N6").Value = -811
$ws.Range("H26").Value = 410
$ws.Range("I26").Value = 200
$ws.Range("J26").Value = 620
$ws.Range("K26").Value = 600
$ws.Range("L26").Value = 1860
$ws.Range("M26").Value = -312
$ws.Range("N26").Value = -2436
$ws.Range("H107").Value = 1023.76
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 1104.0869
$ws.Range("K107").Value = 300
$ws.Range("L107").Value = 3312.2607
$ws.Range("M107").Value = 1620
$ws.Range("N107").Value = -7152.2607
$ws.Range("H124").Value = 5299.8184
$ws.Range("I124").Value = 500
$ws.Range("J124").Value = 5779.8
$ws.Range("K124").Value = 1500
$ws.Range("L124").Value = 17339.4
$ws.Range("M124").Value = 3410
$ws.Range("N124").Value = -27159.4
$ws.Range("H126").Value = 2531.2307
$ws.Range("J126").Value = 2658
$ws.Range("L126").Value = 7974
$ws.Range("N126").Value = -17854
$ws.Range("H130").Value = 6678.3335
$ws.Range("I130").Value = 1628
$ws.Range("J130").Value = 10285.714
$ws.Range("K130").Value = 4884
$ws.Range("L130").Value = 30857.142
$ws.Range("M130").Value = 136
$ws.Range("N130").Value = -40897.142
$ws.Range("H132").Value = 1069.1428
$ws.Range("I132").Value = 830.6667
$ws.Range("K132").Value = 7476.0003
$ws.Range("M132").Value = -4946.0003
$ws.Range("H135").Value = 3109.875
$ws.Range("I135").Value = 434.8
$ws.Range("K135").Value = 3913.2
$ws.Range("M135").Value = -1378.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 3422.5715
$ws.Range("I107").Value = 4340
$ws.Range("K107").Value = 4340
$ws.Range("M107").Value = -2420
$ws.Range("H138").Value = 56874.5
$ws.Range("J138").Value = 56874.5
$ws.Range("L138").Value = 56874.5
$ws.Range("N138").Value = -67154.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 5245.3335
$ws.Range("I122").Value = 5686.9375
$ws.Range("J122").Value = 4829.706
$ws.Range("K122").Value = 17060.8125
$ws.Range("L122").Value = 14489.118
$ws.Range("M122").Value = -14610.8125
$ws.Range("N122").Value = -19389.118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H107").Value = 860.2222
$ws.Range("I107").Value = 930.875
$ws.Range("J107").Value = 295
$ws.Range("K107").Value = 2792.625
$ws.Range("L107").Value = 885
$ws.Range("M107").Value = -872.625
$ws.Range("N107").Value = -4725
$ws.Range("H138").Value = 64656.668
$ws.Range("J138").Value = 64656.668
$ws.Range("L138").Value = 64656.668
$ws.Range("N138").Value = -74936.66800000001
